$wb = $excel.ActiveWorkbook

# --- Teacher Credentials sheet: fix typo'd passcodes (123465 -> 123456) ---
$wsTeacherCred = $wb.Worksheets.Item("Teacher Credentials")
$wsTeacherCred.Range("B2").Value = 123456
$wsTeacherCred.Range("B3").Value = 123456
$wsTeacherCred.Range("B4").Value = 123456
$wsTeacherCred.Range("B5").Value = 123456
$wsTeacherCred.Range("B6").Value = 123456

# Move the selection/cursor on this sheet and make it the active tab
$wsTeacherCred.Activate()
$wsTeacherCred.Range("C5").Select()

# --- TeacherLoginPageCorrectCredenti sheet: move selection to E21 ---
$wsTeacherLoginCorrect = $wb.Worksheets.Item("TeacherLoginPageCorrectCredenti")
$wsTeacherLoginCorrect.Activate()
$wsTeacherLoginCorrect.Range("E21").Select()

# Leave the "Teacher Credentials" tab as the active one when done
$wsTeacherCred.Activate()
